# Daily attendance processing - 2026-01-19 05:48:03
# Normalizes the "Recorded By" column (G) so the System entry is listed
# after the human recorder's e-mail instead of before it, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#
# Rows are walked top to bottom. Whenever two CONSECUTIVE rows both still
# carry the old "System, <email>" value, that matched pair is normalised
# and the walk resumes after it. A row left without a partner (i.e. the
# odd one out at the end of a run of identical values) is skipped, since
# it has already been normalised by an earlier pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$col = 7  # column G - "Recorded By"

# Snapshot the column first so the pairing logic isn't thrown off by the
# edits being made as we go.
$values = @()
for ($r = 1; $r -le $lastRow; $r++) {
    $values += $ws.Cells.Item($r, $col).Value2
}

$r = 1
while ($r -lt $lastRow) {
    if ($values[$r - 1] -eq $oldValue -and $values[$r] -eq $oldValue) {
        $ws.Cells.Item($r, $col).Value = $newValue
        $ws.Cells.Item($r + 1, $col).Value = $newValue
        $r += 2
    }
    else {
        $r += 1
    }
}
